$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# wdCellAlignVerticalBottom = 3
$wdCellAlignVerticalBottom = 3

$rowIndexes = @(1, 2, 5)

foreach ($rowIdx in $rowIndexes) {
    $row = $t.Rows.Item($rowIdx)
    $cellCount = $row.Cells.Count
    for ($c = 1; $c -le $cellCount; $c++) {
        $row.Cells.Item($c).VerticalAlignment = $wdCellAlignVerticalBottom
    }
}
